$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings are preserved as literal text
$dCol = $ws.Range("D2:D51")
$dCol.NumberFormat = "@"

$ws.Range("D2").Value = "33.970.44"
$ws.Range("E2").Value = "  +10.50%  "

$ws.Range("D3").Value = "1.813.38"
$ws.Range("E3").Value = "  +7.34%  "

$ws.Range("D5").Value = "227.96"
$ws.Range("E5").Value = "  +2.96%  "

$ws.Range("D6").Value = "0.540"
$ws.Range("E6").Value = "  +3.49%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "30.88"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").Value = "45.85"
$ws.Range("E9").Value = "  +3.59%  "

$ws.Range("D10").Value = "0.278"
$ws.Range("E10").Value = "  +4.82%  "

$ws.Range("E11").Value = "  +6.70%  "

$ws.Range("E12").Value = "  +2.34%  "

$ws.Range("D13").Value = "2.074.18"
$ws.Range("E13").Value = "  +7.30%  "

$ws.Range("D14").Value = "1.821.05"
$ws.Range("E14").Value = "  +7.69%  "

$ws.Range("D15").Value = "0.638"
$ws.Range("E15").Value = "  +2.16%  "

$ws.Range("D16").Value = "33.984.33"
$ws.Range("E16").Value = "  +10.40%  "

$ws.Range("D17").Value = "10.15"
$ws.Range("E17").Value = "  -4.65%  "

$ws.Range("E18").Value = "  +6.69%  "

$ws.Range("D19").Value = "68.93"
$ws.Range("E19").Value = "  +3.60%  "

$ws.Range("E20").Value = "  +3.24%  "

$ws.Range("D21").Value = "0.0₃0741"
$ws.Range("E21").Value = "  +3.49%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "10.37"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("E25").Value = "  +1.04%  "

$ws.Range("D26").Value = "158.37"
$ws.Range("E26").Value = "  +0.65%  "

$ws.Range("E27").Value = "  +3.91%  "

$ws.Range("E28").Value = "  +3.24%  "

$ws.Range("D29").Value = "7.04"
$ws.Range("E29").Value = "  +4.71%  "

$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("E31").Value = "  +8.83%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0508"
$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  +5.13%  "

$ws.Range("D34").Value = "3.51"
$ws.Range("E34").Value = "  +6.40%  "

$ws.Range("D35").Value = "1.547.57"
$ws.Range("E35").Value = "  +2.07%  "

$ws.Range("D36").Value = "1.81"
$ws.Range("E36").Value = "  +3.50%  "

$ws.Range("E37").Value = "  +3.20%  "

$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").Value = "83.93"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0186"
$ws.Range("E39").Value = "  +3.80%  "

$ws.Range("D40").Value = "0.615"
$ws.Range("E40").Value = "  +5.33%  "

$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  +3.22%  "

$ws.Range("D42").Value = "2.34"
$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("D43").Value = "0.900"
$ws.Range("E43").Value = "  +5.87%  "

$ws.Range("D44").Value = "2.10"
$ws.Range("E44").Value = "  +4.66%  "

$ws.Range("D45").Value = "0.0525"

$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("D47").Value = "1.961.07"
$ws.Range("E47").Value = "  +7.17%  "

$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("D49").Value = "5.65"
$ws.Range("E49").Value = "  +3.67%  "

$ws.Range("D50").Value = "51.86"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "11.09"
$ws.Range("E51").Value = "  +15.40%  "

# Restore default style on column D (removes the text number-format override while keeping values as text)
$dCol.Style = "Normal"

